$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "Feltrinelli Daniele"
$ws.Range("B59").Value = "Stefano Tita | Clitoriders"
$ws.Range("C59").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("D59").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E59").Value = "Luca Frasca | Clitoriders"
$ws.Range("F59").Value = "Davide  Bazzano  | iMontagna"
